$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) The former row 2 ("Omegam_H0 (pour LambdaCDM)" results) moves to
#    row 3. Re-create it there (values identical, formulas re-pointed
#    to the new row).
# ------------------------------------------------------------------
$ws.Range("A3").Value = "Omegam_H0 (pour LambdaCDM)"
$ws.Range("B3").Value = 73.2999999999995
$ws.Range("C3").Value = 1.7
$ws.Range("D3").Formula = "=ABS(B3-(69.5999999999997))"
$ws.Range("E3").Formula = "=33.4000000000001/100"
$ws.Range("F3").Formula = "=ABS(E3-(35.0000000000001/100))"
$ws.Range("G3").Formula = "=ABS(E3-(32.5000000000001/100))"
$ws.Range("K3").Value = 1523.01688482443
$ws.Range("L3").Value = 2
$ws.Range("M3").Formula = "=K3+2*L3"

# Wipe the old row 2 content so we can refill it with the new results.
$ws.Range("A2:M2").ClearContents()

# ------------------------------------------------------------------
# 2) New row 4 ("Omegam_Omegal(pour LambdaCDM)" results). Populated
#    before row 2 below so that its text is registered first in the
#    shared-string table (matching the expected string order).
# ------------------------------------------------------------------
$ws.Range("A4").Value = "Omegam_Omegal(pour LambdaCDM)"
$ws.Range("E4").Value = 0.33
$ws.Range("F4").Formula = "=ABS(E4-0.35)"
$ws.Range("G4").Formula = "=ABS(E4-0.306)"
$ws.Range("H4").Value = 0.66200000000000003
$ws.Range("I4").Formula = "=ABS(H4-0.694)"
$ws.Range("J4").Formula = "=ABS(H4-0.633)"
$ws.Range("K4").Value = 1523.01055414283
$ws.Range("L4").Value = 2
$ws.Range("M4").Formula = "=K4+2*L4"

# ------------------------------------------------------------------
# 3) New row 2 ("Résultats attendus pour Flat LCDM (Brout et al. 2022)")
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Résultats attendus pour Flat LCDM (Brout et al. 2022)"
$ws.Range("B2").Value = 73.599999999999994
$ws.Range("C2").Value = 1.1000000000000001
$ws.Range("D2").Value = 1.1000000000000001
$ws.Range("E2").Value = 0.33400000000000002
$ws.Range("F2").Value = 0.018
$ws.Range("G2").Value = 0.018
$ws.Range("H2").Value = 0.66600000000000004
$ws.Range("I2").Value = 0.018
$ws.Range("J2").Value = 0.018

# ------------------------------------------------------------------
# 4) Column A is widened to fit the new longer labels.
# ------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 49.14

# ------------------------------------------------------------------
# 5) Move the threaded comments that used to sit on row 2 down to row 3
#    (cell content moved, so the annotations follow).
# ------------------------------------------------------------------
$oldC2 = $ws.Range("C2").CommentThreaded
$oldC2Text = $oldC2.Text()
$oldC2.Delete()

$oldG2 = $ws.Range("G2").CommentThreaded
$oldG2Text = $oldG2.Text()
$oldG2.Delete()

$ws.Range("C3").AddCommentThreaded($oldC2Text) | Out-Null
$ws.Range("G3").AddCommentThreaded($oldG2Text) | Out-Null

# ------------------------------------------------------------------
# 6) Two new threaded comments for the newly added row 4 results.
# ------------------------------------------------------------------
$ws.Range("F4").AddCommentThreaded("(sans doute minimisé car range trop petit (s'arrête à 0.35))") | Out-Null
$ws.Range("K4").AddCommentThreaded("Min de chi2 très similaire à celui pour Omegam_H0 : pourrait être utilisé comme référence pour la suite pour voir si les modèles avec les G(z) améliorent ou non les fits") | Out-Null

# ------------------------------------------------------------------
# 7) Selection moved to N4 as the active cell.
# ------------------------------------------------------------------
$ws.Range("N4").Select()

Write-Output "edit applied"
